$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the data table (tone-frame blocks per day, D1-D5) ---
$ws.Range("A1").Value = 'file_path'
$ws.Range("B1").Value = 'start_frame'
$ws.Range("C1").Value = 'stop_frame'
$ws.Range("D1").Value = 'individual_subj'
$ws.Range("E1").Value = 'all_subj'

$ws.Range("A2").Value = ''
$ws.Range("B2").Value = ''
$ws.Range("C2").Value = ''
$ws.Range("D2").Value = ''
$ws.Range("E2").Value = ''

$ws.Range("A3").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2.1.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B3").Value = 41000
$ws.Range("C3").Value = 79051
$ws.Range("D3").Value = 1.1
$ws.Range("E3").Value = '1.1_1.2'

$ws.Range("A4").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2.1.fixed.2_subj.round_1.id_corrected.h5'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 39500
$ws.Range("D4").Value = '1.1_1.2'
$ws.Range("E4").Value = '1.1_1.2'

$ws.Range("A5").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2/20230617_115521_standard_comp_to_omission_D1_subj_1-1_and_1-2.3.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 38957
$ws.Range("D5").Value = 1.2
$ws.Range("E5").Value = '1.1_1.2'

$ws.Range("A6").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1.1.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B6").Value = 32792
$ws.Range("C6").Value = 68495
$ws.Range("D6").Value = 1.4
$ws.Range("E6").Value = '1.1_1.4'

$ws.Range("A7").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1.1.fixed.2_subj.round_1.id_corrected.h5'
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 32316
$ws.Range("D7").Value = '1.1_1.4'
$ws.Range("E7").Value = '1.1_1.4'

$ws.Range("A8").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1/20230618_100636_standard_comp_to_omission_D2_subj_1-4_and_1-1.2.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B8").Value = 32792
$ws.Range("C8").Value = 68495
$ws.Range("D8").Value = 1.1
$ws.Range("E8").Value = '1.1_1.4'

$ws.Range("A9").Value = ''
$ws.Range("B9").Value = ''
$ws.Range("C9").Value = ''
$ws.Range("D9").Value = ''
$ws.Range("E9").Value = ''

$ws.Range("A10").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4.3.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B10").Value = 32500
$ws.Range("C10").Value = 66320
$ws.Range("D10").Value = 1.2
$ws.Range("E10").Value = '1.2_1.4'

$ws.Range("A11").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4.3.fixed.2_subj.round_1.id_corrected.h5'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 30000
$ws.Range("D11").Value = '1.2_1.4'
$ws.Range("E11").Value = '1.2_1.4'

$ws.Range("A12").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4/20230619_115321_standard_comp_to_omission_D3_subj_1-2_and_1-4.4.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B12").Value = 32500
$ws.Range("C12").Value = 66006
$ws.Range("D12").Value = 1.4
$ws.Range("E12").Value = '1.2_1.4'

$ws.Range("A13").Value = ''
$ws.Range("B13").Value = ''
$ws.Range("C13").Value = ''
$ws.Range("D13").Value = ''
$ws.Range("E13").Value = ''

$ws.Range("A14").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1.1.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B14").Value = 32860
$ws.Range("C14").Value = 68288
$ws.Range("D14").Value = 1.2
$ws.Range("E14").Value = '1.1_1.2'

$ws.Range("A15").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1.1.fixed.2_subj.round_1.id_corrected.h5'
$ws.Range("B15").Value = 2027
$ws.Range("C15").Value = 32240
$ws.Range("D15").Value = '1.1_1.2'
$ws.Range("E15").Value = '1.1_1.2'

$ws.Range("A16").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1/20230620_114347_standard_comp_to_omission_D4_subj_1-2_and_1-1.2.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B16").Value = 32860
$ws.Range("C16").Value = 68288
$ws.Range("D16").Value = 1.1
$ws.Range("E16").Value = '1.1_1.2'

$ws.Range("A17").Value = ''
$ws.Range("B17").Value = ''
$ws.Range("C17").Value = ''
$ws.Range("D17").Value = ''
$ws.Range("E17").Value = ''

$ws.Range("A18").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2.1.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B18").Value = 35000
$ws.Range("C18").Value = 69768
$ws.Range("D18").Value = 1.4
$ws.Range("E18").Value = '1.2_1.4'

$ws.Range("A19").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2.1.fixed.2_subj.round_1.id_corrected.h5'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 33500
$ws.Range("D19").Value = '1.2_1.4'
$ws.Range("E19").Value = '1.2_1.4'

$ws.Range("A20").Value = '/scratch/back_up/reward_competition_extention/proc/id_corrected/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2/20230621_111240_standard_comp_to_omission_D5_subj_1-4_and_1-2.2.fixed.1_subj.round_1.id_corrected.h5'
$ws.Range("B20").Value = 34500
$ws.Range("C20").Value = 69684
$ws.Range("D20").Value = 1.2
$ws.Range("E20").Value = '1.2_1.4'

# --- Extend trailing blank placeholder rows down to row 33 ---
$ws.Range("A26").Value = ''
$ws.Range("A27").Value = ''
$ws.Range("A28").Value = ''
$ws.Range("A29").Value = ''
$ws.Range("A30").Value = ''
$ws.Range("A31").Value = ''
$ws.Range("A32").Value = ''
$ws.Range("A33").Value = ''

# --- Update selection to match the authored edit ---
$ws.Range("J16").Select()
